# This script updates the values in the "ldaqda_pc15_perfs" workbook to
# reflect results obtained after adding a validation set method for LDA / QDA.
# It touches three worksheets:
#   1. "Test errors"     - summary error rates for LDA / QDA
#   2. "LDA Conf. Mat."  - LDA confusion matrix
#   3. "QDA Conf. Mat."  - QDA confusion matrix

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet 1: "Test errors"
# ---------------------------------------------------------------------------
$wsErrors = $wb.Worksheets.Item("Test errors")

$wsErrors.Range("B2").Value = 0.17592592592592593
$wsErrors.Range("C2").Value = 0.07790994062208381
$wsErrors.Range("D2").Value = 0.33796296296296297
$wsErrors.Range("E2").Value = 0.05670115145331431

# ---------------------------------------------------------------------------
# Sheet 2: "LDA Conf. Mat."
# ---------------------------------------------------------------------------
$wsLda = $wb.Worksheets.Item("LDA Conf. Mat.")

$wsLda.Range("C3").Value = 34.0
$wsLda.Range("F3").Value = 2.0
$wsLda.Range("G3").Value = 0.0

$wsLda.Range("E4").Value = 1.0
$wsLda.Range("F4").Value = 4.0

$wsLda.Range("D5").Value = 1.0
$wsLda.Range("E5").Value = 29.0
$wsLda.Range("F5").Value = 5.0

$wsLda.Range("D6").Value = 3.0
$wsLda.Range("E6").Value = 5.0

# ---------------------------------------------------------------------------
# Sheet 3: "QDA Conf. Mat."
# ---------------------------------------------------------------------------
$wsQda = $wb.Worksheets.Item("QDA Conf. Mat.")

$wsQda.Range("B2").Value = 19.0
$wsQda.Range("G2").Value = 16.0

$wsQda.Range("C3").Value = 30.0
$wsQda.Range("D3").Value = 1.0

$wsQda.Range("B4").Value = 1.0
$wsQda.Range("D4").Value = 23.0
$wsQda.Range("E4").Value = 1.0
$wsQda.Range("F4").Value = 7.0
$wsQda.Range("G4").Value = 4.0

$wsQda.Range("B5").Value = 0.0
$wsQda.Range("E5").Value = 21.0
$wsQda.Range("F5").Value = 6.0
$wsQda.Range("G5").Value = 8.0

$wsQda.Range("D6").Value = 4.0
$wsQda.Range("E6").Value = 7.0
$wsQda.Range("F6").Value = 24.0
$wsQda.Range("G6").Value = 1.0

$wsQda.Range("B7").Value = 6.0
$wsQda.Range("C7").Value = 1.0
$wsQda.Range("E7").Value = 1.0
$wsQda.Range("F7").Value = 1.0
$wsQda.Range("G7").Value = 26.0
